$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Y"
$ws.Range("D3").Value = "Y"
$ws.Range("D4").Value = "Y"
$ws.Range("D5").Value = "Y"

$ws.Range("D6").Select()
